$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the "Date" metadata row (A8 = "Date") with the new generation timestamp
$ws.Range("B8").Value = "2025-08-20T17:48:34+01:00"

# Fill in the "Description" metadata row (A12 = "Description") which was previously blank
$ws.Range("B12").Value = "Code system for fertility status indicators"
